$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 3, 4) are rotated: row2 -> row4, row3 -> row2, row4 -> row3.
# Additionally, the Ost (Q) / Nord (R) coordinates are rounded to whole numbers,
# and the Starttid (Z) / Sluttid (AB) columns are cleared for every row.

# Capture the "before" values for the columns that actually change between rows.
$A2 = $ws.Range("A2").Value()
$B2 = $ws.Range("B2").Value()
$E2 = $ws.Range("E2").Value()
$F2 = $ws.Range("F2").Value()
$G2 = $ws.Range("G2").Value()
$H2 = $ws.Range("H2").Value()
$P2 = $ws.Range("P2").Value()
$Q2 = $ws.Range("Q2").Value()
$R2 = $ws.Range("R2").Value()
$S2 = $ws.Range("S2").Value()
$AC2 = $ws.Range("AC2").Value()

$A3 = $ws.Range("A3").Value()
$B3 = $ws.Range("B3").Value()
$E3 = $ws.Range("E3").Value()
$F3 = $ws.Range("F3").Value()
$G3 = $ws.Range("G3").Value()
$H3 = $ws.Range("H3").Value()
$P3 = $ws.Range("P3").Value()
$Q3 = $ws.Range("Q3").Value()
$R3 = $ws.Range("R3").Value()
$S3 = $ws.Range("S3").Value()
$AC3 = $ws.Range("AC3").Value()

$A4 = $ws.Range("A4").Value()
$B4 = $ws.Range("B4").Value()
$E4 = $ws.Range("E4").Value()
$F4 = $ws.Range("F4").Value()
$G4 = $ws.Range("G4").Value()
$H4 = $ws.Range("H4").Value()
$P4 = $ws.Range("P4").Value()
$Q4 = $ws.Range("Q4").Value()
$R4 = $ws.Range("R4").Value()
$S4 = $ws.Range("S4").Value()
$AC4 = $ws.Range("AC4").Value()

# New row 2 <= old row 3 (with Q/R rounded)
$ws.Range("A2").Value() = $A3
$ws.Range("B2").Value() = $B3
$ws.Range("E2").Value() = $E3
$ws.Range("F2").Value() = $F3
$ws.Range("G2").Value() = $G3
$ws.Range("H2").Value() = $H3
$ws.Range("P2").Value() = $P3
$ws.Range("Q2").Value() = [math]::Round($Q3)
$ws.Range("R2").Value() = [math]::Round($R3)
$ws.Range("S2").Value() = $S3
$ws.Range("AC2").Value() = $AC3

# New row 3 <= old row 4 (with Q/R rounded)
$ws.Range("A3").Value() = $A4
$ws.Range("B3").Value() = $B4
$ws.Range("E3").Value() = $E4
$ws.Range("F3").Value() = $F4
$ws.Range("G3").Value() = $G4
$ws.Range("H3").Value() = $H4
$ws.Range("P3").Value() = $P4
$ws.Range("Q3").Value() = [math]::Round($Q4)
$ws.Range("R3").Value() = [math]::Round($R4)
$ws.Range("S3").Value() = $S4
$ws.Range("AC3").Value() = $AC4

# New row 4 <= old row 2 (with Q/R rounded)
$ws.Range("A4").Value() = $A2
$ws.Range("B4").Value() = $B2
$ws.Range("E4").Value() = $E2
$ws.Range("F4").Value() = $F2
$ws.Range("G4").Value() = $G2
$ws.Range("H4").Value() = $H2
$ws.Range("P4").Value() = $P2
$ws.Range("Q4").Value() = [math]::Round($Q2)
$ws.Range("R4").Value() = [math]::Round($R2)
$ws.Range("S4").Value() = $S2
$ws.Range("AC4").Value() = $AC2

# Remove Starttid (Z) and Sluttid (AB) for all three data rows.
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
